# Applies the feedback-driven rewording of the Synthesis Project Plan:
# de-personalizes first-person language ("I"/"my") into third-person
# references to "the individual" / "he" / passive voice, per the
# 19 May 2022 review feedback.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        throw "Find/Replace target not found: $find"
    }
}

# 1. Introduction paragraph.
Replace-Text `
    "I will be working on during weeks 11-16. This is an individual assignment where I can show my current proficiency of all seven (7) learning outcomes (LOs)." `
    "on which the work will be done during weeks 11-16. This is an individual assignment where the current proficiency of all seven (7) learning outcomes (LOs) can be showed."

# 2. Team paragraph ("I am the representant" -> "The individual is the representant", "I" -> "He"/"he").
Replace-Text `
    "I am the representant. I will be in contact with the tutor and will be responsible for organizing meetings. During the meetings I will be taking minutes. " `
    "The individual is the representant. He will be in contact with the tutor and will be responsible for organizing meetings. During the meetings he will be taking minutes. "

# 3. Deliverables paragraph (active -> passive).
Replace-Text `
    "By the end of the allotted six weeks, I will deliver a fully working application and website that addresses the main issues." `
    "By the end of the allotted six weeks, a fully working application and website that addresses the main issues will be delivered."

# 4. "The company [...] is having the following problems ..."
Replace-Text `
    " is having the following problems and we will make a software solution able to solve them:" `
    " is having the following problems and this software solution will be able to solve them:"

# 5. Non-deliverables intro paragraph.
Replace-Text `
    "These non-deliverables are what the company could expect me to deliver but I will not, as I will focus on providing the grounds for the required actions as well providing a test plan and a test report." `
    "These non-deliverables are what the company could expect to be delivered but they will not, as the focus will be mainly on providing the grounds for the required actions as well providing a test plan and a test report."

# 6. Written manual non-deliverable bullet.
Replace-Text `
    "I will not deliver any written manual for our application." `
    "A written manual for the application will not be delivered."

# 7. Constraints: programming language paragraph.
Replace-Text `
    "I am going to use C# for the desktop part of the application. It has a defined style and I do not plan to deviate too much from that. " `
    "The C# programming language will be used for the desktop part of the application. It has a defined style and there will be not much derivation from that. "

# 8. Constraints: IDE paragraph.
Replace-Text `
    "for this reason I will be using Microsoft Visual Studio." `
    "for this reason the individual will be using Microsoft Visual Studio."

# 9. Constraints: pre-defined time paragraph.
Replace-Text `
    "I have pre-defined amount of time to complete the project which will limit the additional features I can provide. " `
    "The individual has pre-defined the amount of time required to complete the project which will limit the additional features that can be provided. "

# 10. Risks paragraph.
Replace-Text `
    "The only risk I foresee is getting tangled up in the extra requirements as I want to provide more than the bare minimum." `
    "The only risk foreseen is getting tangled up in the extra requirements as the individual wants to provide more than the bare minimum."

Write-Output "All 10 replacements applied successfully."
